# Trading update: 2026-02-18 10:32:02
$wb = $excel.ActiveWorkbook

# ---- All Trades sheet: append rows 30 and 31 ----
$wsAll = $wb.Worksheets.Item("All Trades")

# Row 30 (Trade #29 - momentum)
$wsAll.Cells.Item(30, 1).Value = 29
$wsAll.Cells.Item(30, 2).Value = "'2026-02-18"
$wsAll.Cells.Item(30, 3).Value = "10:30:12"
$wsAll.Cells.Item(30, 4).Value = "momentum"
$wsAll.Cells.Item(30, 5).Value = "UP"
$wsAll.Cells.Item(30, 6).Value = 0.01
$wsAll.Cells.Item(30, 8).Value = "OPEN"
$wsAll.Cells.Item(30, 9).Value = 0
$wsAll.Cells.Item(30, 10).Value = 0
$wsAll.Cells.Item(30, 11).Value = 100
$wsAll.Cells.Item(30, 13).Value = 0
$wsAll.Cells.Item(30, 14).Value = 0
$wsAll.Cells.Item(30, 15).Value = 0
$wsAll.Cells.Item(30, 16).Value = 0.9
$wsAll.Cells.Item(30, 17).Value = "Upward momentum: 1.020% over 10 samples"
# G30 (Exit Price) and L30 (Exit Reason) stay blank like the row above - copy
# the existing blank cells so the row30/col7/col12 positions exist in the sheet.
$wsAll.Cells.Item(29, 7).Copy($wsAll.Cells.Item(30, 7))
$wsAll.Cells.Item(29, 12).Copy($wsAll.Cells.Item(30, 12))

# Row 31 (Trade #30 - MarketMaking)
$wsAll.Cells.Item(31, 1).Value = 30
$wsAll.Cells.Item(31, 2).Value = "'2026-02-18"
$wsAll.Cells.Item(31, 3).Value = "10:30:12"
$wsAll.Cells.Item(31, 4).Value = "MarketMaking"
$wsAll.Cells.Item(31, 5).Value = "UP"
$wsAll.Cells.Item(31, 6).Value = 0.01
$wsAll.Cells.Item(31, 8).Value = "OPEN"
$wsAll.Cells.Item(31, 9).Value = 0
$wsAll.Cells.Item(31, 10).Value = 0
$wsAll.Cells.Item(31, 11).Value = 100
$wsAll.Cells.Item(31, 13).Value = 0
$wsAll.Cells.Item(31, 14).Value = 0
$wsAll.Cells.Item(31, 15).Value = 0
$wsAll.Cells.Item(31, 16).Value = 0.6
$wsAll.Cells.Item(31, 17).Value = "Normal spread capture: 202 bps"
$wsAll.Cells.Item(29, 7).Copy($wsAll.Cells.Item(31, 7))
$wsAll.Cells.Item(29, 12).Copy($wsAll.Cells.Item(31, 12))

# ---- momentum sheet: append row 4 ----
$wsMom = $wb.Worksheets.Item("momentum")

$wsMom.Cells.Item(4, 1).Value = 29
$wsMom.Cells.Item(4, 2).Value = "'2026-02-18"
$wsMom.Cells.Item(4, 3).Value = "10:30:12"
$wsMom.Cells.Item(4, 4).Value = "momentum"
$wsMom.Cells.Item(4, 5).Value = "UP"
$wsMom.Cells.Item(4, 6).Value = 0.01
$wsMom.Cells.Item(4, 8).Value = "OPEN"
$wsMom.Cells.Item(4, 9).Value = 0
$wsMom.Cells.Item(4, 10).Value = 0
$wsMom.Cells.Item(4, 11).Value = 100
$wsMom.Cells.Item(4, 12).Value = 0
$wsMom.Cells.Item(4, 13).Value = 0
$wsMom.Cells.Item(4, 14).Value = 0.9
$wsMom.Cells.Item(4, 15).Value = "Upward momentum: 1.020% over 10 samples"
$wsMom.Cells.Item(4, 17).Value = 0
# G4 (Exit Price) and P4 (Exit Reason) stay blank like the row above.
$wsMom.Cells.Item(3, 7).Copy($wsMom.Cells.Item(4, 7))
$wsMom.Cells.Item(3, 16).Copy($wsMom.Cells.Item(4, 16))

# ---- MarketMaking sheet: append row 6 ----
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Cells.Item(6, 1).Value = 30
$wsMM.Cells.Item(6, 2).Value = "'2026-02-18"
$wsMM.Cells.Item(6, 3).Value = "10:30:12"
$wsMM.Cells.Item(6, 4).Value = "MarketMaking"
$wsMM.Cells.Item(6, 5).Value = "UP"
$wsMM.Cells.Item(6, 6).Value = 0.01
$wsMM.Cells.Item(6, 8).Value = "OPEN"
$wsMM.Cells.Item(6, 9).Value = 0
$wsMM.Cells.Item(6, 10).Value = 0
$wsMM.Cells.Item(6, 11).Value = 100
$wsMM.Cells.Item(6, 12).Value = 0
$wsMM.Cells.Item(6, 13).Value = 0
$wsMM.Cells.Item(6, 14).Value = 0.6
$wsMM.Cells.Item(6, 15).Value = "Normal spread capture: 202 bps"
$wsMM.Cells.Item(6, 17).Value = 0
# G6 (Exit Price) and P6 (Exit Reason) stay blank like the row above.
$wsMM.Cells.Item(5, 7).Copy($wsMM.Cells.Item(6, 7))
$wsMM.Cells.Item(5, 16).Copy($wsMM.Cells.Item(6, 16))
